$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(18, 9).Value = "sd"
$ws.Cells.Item(18, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(24, 9).Value = "qy"
$ws.Cells.Item(24, 10).Value = "Yes-No-Question"
$ws.Cells.Item(25, 9).Value = "sv"
$ws.Cells.Item(25, 10).Value = "Statement-opinion"
$ws.Cells.Item(28, 9).Value = "sd"
$ws.Cells.Item(28, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(31, 9).Value = "aa"
$ws.Cells.Item(31, 10).Value = "Agree/Accept"
$ws.Cells.Item(35, 9).Value = "sd"
$ws.Cells.Item(35, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(42, 9).Value = "%"
$ws.Cells.Item(42, 10).Value = "Uninterpretable"
$ws.Cells.Item(53, 9).Value = "sd"
$ws.Cells.Item(53, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(63, 9).Value = "aa"
$ws.Cells.Item(63, 10).Value = "Agree/Accept"
$ws.Cells.Item(66, 9).Value = "aa"
$ws.Cells.Item(66, 10).Value = "Agree/Accept"
$ws.Cells.Item(67, 9).Value = "sv"
$ws.Cells.Item(67, 10).Value = "Statement-opinion"
$ws.Cells.Item(72, 9).Value = "aa"
$ws.Cells.Item(72, 10).Value = "Agree/Accept"
$ws.Cells.Item(83, 9).Value = "sv"
$ws.Cells.Item(83, 10).Value = "Statement-opinion"
$ws.Cells.Item(85, 9).Value = "aa"
$ws.Cells.Item(85, 10).Value = "Agree/Accept"
$ws.Cells.Item(103, 9).Value = "sv"
$ws.Cells.Item(103, 10).Value = "Statement-opinion"
$ws.Cells.Item(104, 9).Value = "sd"
$ws.Cells.Item(104, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(109, 9).Value = "sv"
$ws.Cells.Item(109, 10).Value = "Statement-opinion"
$ws.Cells.Item(110, 9).Value = "%"
$ws.Cells.Item(110, 10).Value = "Uninterpretable"
$ws.Cells.Item(127, 9).Value = "aa"
$ws.Cells.Item(127, 10).Value = "Agree/Accept"
$ws.Cells.Item(131, 9).Value = "sd"
$ws.Cells.Item(131, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(133, 9).Value = "sd"
$ws.Cells.Item(133, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(141, 9).Value = "sv"
$ws.Cells.Item(141, 10).Value = "Statement-opinion"
$ws.Cells.Item(149, 9).Value = "sd"
$ws.Cells.Item(149, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(151, 9).Value = "sv"
$ws.Cells.Item(151, 10).Value = "Statement-opinion"
$ws.Cells.Item(154, 9).Value = "%"
$ws.Cells.Item(154, 10).Value = "Uninterpretable"
$ws.Cells.Item(157, 9).Value = "sv"
$ws.Cells.Item(157, 10).Value = "Statement-opinion"
$ws.Cells.Item(158, 9).Value = "sd"
$ws.Cells.Item(158, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(167, 9).Value = "sd"
$ws.Cells.Item(167, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(175, 9).Value = "%"
$ws.Cells.Item(175, 10).Value = "Uninterpretable"
$ws.Cells.Item(177, 9).Value = "sv"
$ws.Cells.Item(177, 10).Value = "Statement-opinion"
$ws.Cells.Item(185, 9).Value = "sd"
$ws.Cells.Item(185, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(187, 9).Value = "ba"
$ws.Cells.Item(187, 10).Value = "Appreciation"
$ws.Cells.Item(189, 9).Value = "%"
$ws.Cells.Item(189, 10).Value = "Uninterpretable"
$ws.Cells.Item(195, 9).Value = "sd"
$ws.Cells.Item(195, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(196, 9).Value = "aa"
$ws.Cells.Item(196, 10).Value = "Agree/Accept"
$ws.Cells.Item(198, 9).Value = "aa"
$ws.Cells.Item(198, 10).Value = "Agree/Accept"
$ws.Cells.Item(204, 9).Value = "ba"
$ws.Cells.Item(204, 10).Value = "Appreciation"
$ws.Cells.Item(205, 9).Value = "sd"
$ws.Cells.Item(205, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(216, 9).Value = "sd"
$ws.Cells.Item(216, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(223, 9).Value = "sv"
$ws.Cells.Item(223, 10).Value = "Statement-opinion"
$ws.Cells.Item(253, 9).Value = "aa"
$ws.Cells.Item(253, 10).Value = "Agree/Accept"
$ws.Cells.Item(254, 9).Value = "sd"
$ws.Cells.Item(254, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(256, 9).Value = "sd"
$ws.Cells.Item(256, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(260, 9).Value = "sd"
$ws.Cells.Item(260, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(265, 9).Value = "sd"
$ws.Cells.Item(265, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(266, 9).Value = "sd"
$ws.Cells.Item(266, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(268, 9).Value = "aa"
$ws.Cells.Item(268, 10).Value = "Agree/Accept"
$ws.Cells.Item(270, 9).Value = "aa"
$ws.Cells.Item(270, 10).Value = "Agree/Accept"
$ws.Cells.Item(281, 9).Value = "aa"
$ws.Cells.Item(281, 10).Value = "Agree/Accept"
$ws.Cells.Item(287, 9).Value = "b"
$ws.Cells.Item(287, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(289, 9).Value = "b"
$ws.Cells.Item(289, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(298, 9).Value = "%"
$ws.Cells.Item(298, 10).Value = "Uninterpretable"
$ws.Cells.Item(304, 9).Value = "aa"
$ws.Cells.Item(304, 10).Value = "Agree/Accept"
$ws.Cells.Item(306, 9).Value = "sd"
$ws.Cells.Item(306, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(326, 9).Value = "aa"
$ws.Cells.Item(326, 10).Value = "Agree/Accept"
$ws.Cells.Item(332, 9).Value = "sd"
$ws.Cells.Item(332, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(354, 9).Value = "sd"
$ws.Cells.Item(354, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(356, 9).Value = "sv"
$ws.Cells.Item(356, 10).Value = "Statement-opinion"
$ws.Cells.Item(357, 9).Value = "sd"
$ws.Cells.Item(357, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(362, 9).Value = "sd"
$ws.Cells.Item(362, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(398, 9).Value = "sd"
$ws.Cells.Item(398, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(416, 9).Value = "sd"
$ws.Cells.Item(416, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(422, 9).Value = "sd"
$ws.Cells.Item(422, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(424, 9).Value = "sd"
$ws.Cells.Item(424, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(426, 9).Value = "sd"
$ws.Cells.Item(426, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(436, 9).Value = "sd"
$ws.Cells.Item(436, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(470, 9).Value = "sd"
$ws.Cells.Item(470, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(472, 9).Value = "%"
$ws.Cells.Item(472, 10).Value = "Uninterpretable"
$ws.Cells.Item(479, 9).Value = "b"
$ws.Cells.Item(479, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(488, 9).Value = "sd"
$ws.Cells.Item(488, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(490, 9).Value = "aa"
$ws.Cells.Item(490, 10).Value = "Agree/Accept"
$ws.Cells.Item(493, 9).Value = "aa"
$ws.Cells.Item(493, 10).Value = "Agree/Accept"
$ws.Cells.Item(494, 9).Value = "sd"
$ws.Cells.Item(494, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(522, 9).Value = "sv"
$ws.Cells.Item(522, 10).Value = "Statement-opinion"
$ws.Cells.Item(528, 9).Value = "sd"
$ws.Cells.Item(528, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(544, 9).Value = "aa"
$ws.Cells.Item(544, 10).Value = "Agree/Accept"
$ws.Cells.Item(554, 9).Value = "aa"
$ws.Cells.Item(554, 10).Value = "Agree/Accept"
$ws.Cells.Item(555, 9).Value = "aa"
$ws.Cells.Item(555, 10).Value = "Agree/Accept"
$ws.Cells.Item(559, 9).Value = "aa"
$ws.Cells.Item(559, 10).Value = "Agree/Accept"
$ws.Cells.Item(560, 9).Value = "aa"
$ws.Cells.Item(560, 10).Value = "Agree/Accept"
$ws.Cells.Item(563, 9).Value = "sv"
$ws.Cells.Item(563, 10).Value = "Statement-opinion"
$ws.Cells.Item(572, 9).Value = "sv"
$ws.Cells.Item(572, 10).Value = "Statement-opinion"
$ws.Cells.Item(576, 9).Value = "aa"
$ws.Cells.Item(576, 10).Value = "Agree/Accept"
$ws.Cells.Item(581, 9).Value = "sd"
$ws.Cells.Item(581, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(593, 9).Value = "aa"
$ws.Cells.Item(593, 10).Value = "Agree/Accept"
$ws.Cells.Item(595, 9).Value = "sd"
$ws.Cells.Item(595, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(597, 9).Value = "sd"
$ws.Cells.Item(597, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(599, 9).Value = "sd"
$ws.Cells.Item(599, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(603, 9).Value = "aa"
$ws.Cells.Item(603, 10).Value = "Agree/Accept"
$ws.Cells.Item(615, 9).Value = "sd"
$ws.Cells.Item(615, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(618, 9).Value = "ba"
$ws.Cells.Item(618, 10).Value = "Appreciation"
$ws.Cells.Item(623, 9).Value = "sd"
$ws.Cells.Item(623, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(626, 9).Value = "sv"
$ws.Cells.Item(626, 10).Value = "Statement-opinion"
$ws.Cells.Item(632, 9).Value = "sv"
$ws.Cells.Item(632, 10).Value = "Statement-opinion"
$ws.Cells.Item(639, 9).Value = "sv"
$ws.Cells.Item(639, 10).Value = "Statement-opinion"
$ws.Cells.Item(640, 9).Value = "%"
$ws.Cells.Item(640, 10).Value = "Uninterpretable"
$ws.Cells.Item(642, 9).Value = "sd"
$ws.Cells.Item(642, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(648, 9).Value = "sv"
$ws.Cells.Item(648, 10).Value = "Statement-opinion"
$ws.Cells.Item(659, 9).Value = "sd"
$ws.Cells.Item(659, 10).Value = "Statement-non-opinion"
